$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 'rise'
$ws.Cells.Item(2, 4).Value = 0.7896191477775574
$ws.Cells.Item(3, 3).Value = 'high'
$ws.Cells.Item(3, 4).Value = 0.7600307464599609
$ws.Cells.Item(4, 3).Value = 'demand'
$ws.Cells.Item(4, 4).Value = 0.6293131113052368
$ws.Cells.Item(5, 3).Value = 'uncertainty'
$ws.Cells.Item(5, 4).Value = 0.6166055202484131
$ws.Cells.Item(6, 3).Value = 'mortgage'
$ws.Cells.Item(6, 4).Value = 0.6061130166053772
$ws.Cells.Item(7, 3).Value = 'continue'
$ws.Cells.Item(7, 4).Value = 0.5932634472846985
$ws.Cells.Item(8, 3).Value = 'inflation'
$ws.Cells.Item(8, 4).Value = 0.566853404045105
$ws.Cells.Item(9, 3).Value = 'elevated'
$ws.Cells.Item(9, 4).Value = 0.5618152022361755
$ws.Cells.Item(10, 3).Value = 'commercial'
$ws.Cells.Item(10, 4).Value = 0.5609279870986938
$ws.Cells.Item(11, 3).Value = 'concern'
$ws.Cells.Item(11, 4).Value = 0.5577398538589478
$ws.Cells.Item(12, 3).Value = 'interest'
$ws.Cells.Item(12, 4).Value = 0.566853404045105
$ws.Cells.Item(13, 3).Value = 'rise'
$ws.Cells.Item(13, 4).Value = 0.5616406202316284
$ws.Cells.Item(14, 3).Value = 'pressure'
$ws.Cells.Item(14, 4).Value = 0.5366746783256531
$ws.Cells.Item(15, 3).Value = 'high'
$ws.Cells.Item(15, 4).Value = 0.5099443793296814
$ws.Cells.Item(16, 3).Value = 'employment'
$ws.Cells.Item(16, 4).Value = 0.5018878579139709
$ws.Cells.Item(17, 3).Value = 'concern'
$ws.Cells.Item(17, 4).Value = 0.5006101727485657
$ws.Cells.Item(18, 3).Value = 'contact'
$ws.Cells.Item(18, 4).Value = 0.4986400604248047
$ws.Cells.Item(19, 3).Value = 'continue'
$ws.Cells.Item(19, 4).Value = 0.496562659740448
$ws.Cells.Item(20, 3).Value = 'expectation'
$ws.Cells.Item(20, 4).Value = 0.4741524457931518
$ws.Cells.Item(21, 3).Value = 'increase'
$ws.Cells.Item(21, 4).Value = 0.4714706838130951
$ws.Cells.Item(22, 3).Value = 'automation'
$ws.Cells.Item(22, 4).Value = 0.4898182451725006
$ws.Cells.Item(23, 3).Value = 'laborsave'
$ws.Cells.Item(23, 4).Value = 0.4809595048427582
$ws.Cells.Item(24, 3).Value = 'design'
$ws.Cells.Item(24, 4).Value = 0.4202724099159241
$ws.Cells.Item(25, 3).Value = 'downsizing'
$ws.Cells.Item(25, 4).Value = 0.405096024274826
$ws.Cells.Item(26, 3).Value = 'adoption'
$ws.Cells.Item(26, 4).Value = 0.3957929015159607
$ws.Cells.Item(27, 3).Value = 'endus'
$ws.Cells.Item(27, 4).Value = 0.3899362683296203
$ws.Cells.Item(28, 3).Value = 'interestbearing'
$ws.Cells.Item(28, 4).Value = 0.3792942762374878
$ws.Cells.Item(29, 3).Value = 'technology'
$ws.Cells.Item(29, 4).Value = 0.3787357211112976
$ws.Cells.Item(30, 3).Value = 'composition'
$ws.Cells.Item(30, 4).Value = 0.3747909367084503
$ws.Cells.Item(31, 3).Value = 'tiously'
$ws.Cells.Item(31, 4).Value = 0.3745763599872589
$ws.Cells.Item(32, 3).Value = 'skilled'
$ws.Cells.Item(32, 4).Value = 0.5868028402328491
$ws.Cells.Item(33, 3).Value = 'poach'
$ws.Cells.Item(33, 4).Value = 0.4135763943195343
$ws.Cells.Item(34, 3).Value = 'hourly'
$ws.Cells.Item(34, 4).Value = 0.4074746668338775
$ws.Cells.Item(35, 3).Value = 'phase'
$ws.Cells.Item(35, 4).Value = 0.3961473107337951
$ws.Cells.Item(36, 3).Value = 'tremendous'
$ws.Cells.Item(36, 4).Value = 0.3845434188842773
$ws.Cells.Item(37, 3).Value = 'worker'
$ws.Cells.Item(37, 4).Value = 0.3815480768680572
$ws.Cells.Item(38, 3).Value = 'ture'
$ws.Cells.Item(38, 4).Value = 0.3786315321922302
$ws.Cells.Item(39, 3).Value = 'aforementioned'
$ws.Cells.Item(39, 4).Value = 0.3781601190567016
$ws.Cells.Item(40, 3).Value = 'machinist'
$ws.Cells.Item(40, 4).Value = 0.3765192925930023
$ws.Cells.Item(41, 3).Value = 'refinery'
$ws.Cells.Item(41, 4).Value = 0.3755677342414856
$ws.Cells.Item(42, 3).Value = 'outlook'
$ws.Cells.Item(42, 4).Value = 0.5144216418266296
$ws.Cells.Item(43, 3).Value = 'highly'
$ws.Cells.Item(43, 4).Value = 0.4922049343585968
$ws.Cells.Item(44, 3).Value = 'optimistic'
$ws.Cells.Item(44, 4).Value = 0.4570409953594208
$ws.Cells.Item(45, 3).Value = 'forward'
$ws.Cells.Item(45, 4).Value = 0.4426227807998657
$ws.Cells.Item(46, 3).Value = 'loom'
$ws.Cells.Item(46, 4).Value = 0.4351649582386017
$ws.Cells.Item(47, 3).Value = 'decidedly'
$ws.Cells.Item(47, 4).Value = 0.4260418117046356
$ws.Cells.Item(48, 3).Value = 'navigate'
$ws.Cells.Item(48, 4).Value = 0.4253073930740356
$ws.Cells.Item(49, 3).Value = 'administer'
$ws.Cells.Item(49, 4).Value = 0.4057942032814026
$ws.Cells.Item(50, 3).Value = 'pessimistic'
$ws.Cells.Item(50, 4).Value = 0.4022544622421264
$ws.Cells.Item(51, 3).Value = 'impact'
$ws.Cells.Item(51, 4).Value = 0.3973610401153564
